$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (ID Competição) for data rows 2 through 66 from 65 to 265
$ws.Range("B2:B66").Value = 265
